$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19: "Found more assets and a font to use in the game"
$ws.Range("B19").Value = 45373
$ws.Range("B19").NumberFormat = "d-mmm"
$ws.Range("B19").Font.Size = 11
$ws.Range("B19").Font.Name = "Calibri"
$ws.Range("B19").Font.ThemeColor = 1
$ws.Range("C19").Value = 1
$ws.Range("D19").Value = "Found more assets and a font to use in the game"

# Row 20: "Added assets and font to game. Made game look good" / collision note
$ws.Range("B20").Value = 45374
$ws.Range("B20").NumberFormat = "d-mmm"
$ws.Range("B20").Font.Size = 11
$ws.Range("B20").Font.Name = "Calibri"
$ws.Range("B20").Font.ThemeColor = 1
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = "Added assets and font to game. Made game look good"
$ws.Range("E20").Value = "Collision with missle wont work. Still need to fix the AI then work on the players attacks, finally think about having a second map that shows the player where they have hit more miss to increase playablity like in real battleships. "

$ws.Rows(18).RowHeight = 27
$ws.Rows(20).RowHeight = 28

$ws.Range("E22").Select() | Out-Null
